# Updates cryptos list values/percentages per upstream diff (rows 2-51).
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h).
# D-column assignments are prefixed with a leading apostrophe so that
# Excel's COM layer keeps purely-numeric-looking text (e.g. '600.29')
# stored as text instead of auto-converting it to a numeric value,
# matching the original inlineStr text cells in the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''65.673.06'
$ws.Range('E2').Value = '  -0.31%  '

$ws.Range('D3').Value = '''2.675.09'
$ws.Range('E3').Value = '  -1.16%  '

$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('D5').Value = '''600.29'
$ws.Range('E5').Value = '  -1.58%  '

$ws.Range('D6').Value = '''156.47'
$ws.Range('E6').Value = '  -1.02%  '

$ws.Range('D8').Value = '''0.616'
$ws.Range('E8').Value = '  +4.65%  '

$ws.Range('E9').Value = '  +3.50%  '

$ws.Range('D10').Value = '''5.88'
$ws.Range('E10').Value = '  -2.04%  '

$ws.Range('D11').Value = '''0.399'
$ws.Range('E11').Value = '  -1.07%  '

$ws.Range('E12').Value = '  -0.11%  '

$ws.Range('D13').Value = '''29.24'
$ws.Range('E13').Value = '  -3.92%  '

$ws.Range('D14').Value = '''0.0000196'
$ws.Range('E14').Value = '  -3.61%  '

$ws.Range('D15').Value = '''3.157.00'
$ws.Range('E15').Value = '  -1.14%  '

$ws.Range('D16').Value = '''65.555.87'
$ws.Range('E16').Value = '  -0.33%  '

$ws.Range('D17').Value = '''2.673.42'
$ws.Range('E17').Value = '  -1.33%  '

$ws.Range('D18').Value = '''12.88'
$ws.Range('E18').Value = '  +1.04%  '

$ws.Range('E19').Value = '  -2.17%  '

$ws.Range('E20').Value = '  -1.18%  '

$ws.Range('D21').Value = '''352.61'
$ws.Range('E21').Value = '  -2.07%  '

$ws.Range('E22').Value = '  -0.01%  '

$ws.Range('D23').Value = '''69.77'
$ws.Range('E23').Value = '  -1.67%  '

$ws.Range('D24').Value = '''0.0000112'
$ws.Range('E24').Value = '  +4.73%  '

$ws.Range('D25').Value = '''9.64'
$ws.Range('E25').Value = '  -2.68%  '

$ws.Range('D26').Value = '''1.67'
$ws.Range('E26').Value = '  -0.09%  '

$ws.Range('E27').Value = '  -2.79%  '

$ws.Range('D28').Value = '''1.59'
$ws.Range('E28').Value = '  -5.97%  '

$ws.Range('D29').Value = '''8.04'
$ws.Range('E29').Value = '  -4.56%  '

$ws.Range('E30').Value = '  +0.19%  '

$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = '''2.15'
$ws.Range('E31').Value = '  -3.02%  '

$ws.Range('B32').Value = 'Bittensor'
$ws.Range('C32').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D32').Value = '''529.06'
$ws.Range('E32').Value = '  -2.82%  '

$ws.Range('E33').Value = '  -3.26%  '

$ws.Range('D34').Value = '''5.56'
$ws.Range('E34').Value = '  +1.84%  '

$ws.Range('D35').Value = '''6.51'
$ws.Range('E35').Value = '  -3.16%  '

$ws.Range('E36').Value = '  -2.76%  '

$ws.Range('D37').Value = '''20.59'
$ws.Range('E37').Value = '  -1.75%  '

$ws.Range('E38').Value = '  -0.01%  '

$ws.Range('D39').Value = '''158.29'
$ws.Range('E39').Value = '  -3.19%  '

$ws.Range('E40').Value = '  -2.60%  '

$ws.Range('E41').Value = '  +0.02%  '

$ws.Range('D42').Value = '''163.75'
$ws.Range('E42').Value = '  -5.52%  '

$ws.Range('D43').Value = '''4.13'
$ws.Range('E43').Value = '  -1.87%  '

$ws.Range('D44').Value = '''2.33'
$ws.Range('E44').Value = '  +1.96%  '

$ws.Range('D45').Value = '''0.0609'
$ws.Range('E45').Value = '  -1.58%  '

$ws.Range('D46').Value = '''22.79'
$ws.Range('E46').Value = '  -3.85%  '

$ws.Range('E47').Value = '  -3.37%  '

$ws.Range('D48').Value = '''0.639'
$ws.Range('E48').Value = '  -2.69%  '

$ws.Range('D49').Value = '''0.0₆0262'
$ws.Range('E49').Value = '  +14.02%  '

$ws.Range('D50').Value = '''20.24'
$ws.Range('E50').Value = '  -4.23%  '

$ws.Range('D51').Value = '''0.0997'
$ws.Range('E51').Value = '  +0.35%  '
